# Update "想去人数" (interest count) values in the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 592
$ws1.Range("F7").Value = 14812
$ws1.Range("F12").Value = 8688
$ws1.Range("F13").Value = 328
$ws1.Range("F25").Value = 1075
$ws1.Range("F30").Value = 30
$ws1.Range("F38").Value = 5342

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 592
$ws4.Range("F7").Value = 14812
$ws4.Range("F12").Value = 8688
$ws4.Range("F13").Value = 328
$ws4.Range("F26").Value = 1075
$ws4.Range("F31").Value = 30
$ws4.Range("F41").Value = 5342
